# Update du plan d'adressage Cisco : le reseau "imprimante/tel" passe d'un
# /28 (192.168.5.0/28) a un /27 (192.168.5.0/27), ce qui modifie le masque,
# la derniere IP utilisable, l'adresse de broadcast et le nombre max de
# machines pour cette plage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F20").Value = "192.168.5.0/27"
$ws.Range("F21").Value = "255.255.255.224"
$ws.Range("F23").Value = "192.168.5.27"
$ws.Range("F24").Value = "192.168.5.28"
$ws.Range("F25").Value = "255.255.255.224"
$ws.Range("F26").Value = 28

# Mise a jour cosmetique de la vue (zoom + cellule selectionnee), pour
# coller a l'etat enregistre par l'auteur.
$ws.Activate()
$ws.Range("F23").Select()
$excel.ActiveWindow.Zoom = 69
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
